$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.777.41'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '''1.546.69'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''204.48'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("E6").Value = '  -1.69%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = '''21.31'
$ws.Range("E9").Value = '  -4.41%  '
$ws.Range("D10").Value = '''0.0581'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").Value = '''1.765.36'
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("D13").Value = '''1.543.95'
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '''26.772.97'
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").Value = '''60.92'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '''213.39'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = '''7.25'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").Value = '''1.99'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").Value = '''152.69'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  -2.98%  '
$ws.Range("D27").Value = '''14.80'
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  -2.51%  '
$ws.Range("D30").Value = '''0.0461'
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -3.45%  '
$ws.Range("D32").Value = '''3.17'
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("D33").Value = '''1.352.18'
$ws.Range("E33").Value = '  -3.83%  '
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("E37").Value = '  -3.33%  '
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").Value = '''0.799'
$ws.Range("E40").Value = '  -2.50%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = '''0.992'
$ws.Range("E42").Value = '  -0.98%  '
$ws.Range("D43").Value = '''5.51'
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("E45").Value = '  -2.59%  '
$ws.Range("D46").Value = '''62.79'
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").Value = '''1.680.07'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").Value = '''85.75'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("E50").Value = '  +2.33%  '
$ws.Range("D51").Value = '''0.0₇0975'
$ws.Range("E51").Value = '  -1.40%  '
